# Scheduled market-data refresh: update Leve profit calculations (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets with freshly fetched
# average-price data. Generated from the upstream commit diff.
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 9
$ws.Range("H9").Value = 219.8
$ws.Range("I9").Value = 216.66667
$ws.Range("K9").Value = 216.66667
$ws.Range("M9").Value = -47.66667000000001
# row 12
$ws.Range("H12").Value = 298
$ws.Range("I12").Value = 298
$ws.Range("K12").Value = 298
$ws.Range("M12").Value = -128
# row 40
$ws.Range("H40").Value = 3469.8
$ws.Range("J40").Value = 3697.6
$ws.Range("L40").Value = 3697.6
$ws.Range("N40").Value = -4047.6
# row 48
$ws.Range("H48").Value = 1663
$ws.Range("J48").Value = 2000
$ws.Range("L48").Value = 6000
$ws.Range("N48").Value = -6584
# row 56
$ws.Range("H56").Value = 1663
$ws.Range("J56").Value = 2000
$ws.Range("L56").Value = 6000
$ws.Range("N56").Value = -7068
# row 113
$ws.Range("H113").Value = 100003080
$ws.Range("I113").Value = 33336468
$ws.Range("K113").Value = 33336468
$ws.Range("M113").Value = -33333214

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 57
$ws.Range("H57").Value = 19912
$ws.Range("I57").Value = 19912
$ws.Range("K57").Value = 19912
$ws.Range("M57").Value = -19428
# row 61
$ws.Range("H61").Value = 71586350
$ws.Range("I61").Value = 125015000
$ws.Range("K61").Value = 125015000
$ws.Range("M61").Value = -125014788
# row 136
$ws.Range("H136").Value = 71586350
$ws.Range("I136").Value = 125015000
$ws.Range("K136").Value = 375045000
$ws.Range("M136").Value = -375042450

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 108
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
# row 132
$ws.Range("H132").Value = 94450
$ws.Range("J132").Value = 89666.664
$ws.Range("L132").Value = 89666.664
$ws.Range("N132").Value = -99786.664
# row 134
$ws.Range("H134").Value = 65225.875
$ws.Range("I134").Value = 2200.6667
$ws.Range("K134").Value = 6602.000100000001
$ws.Range("M134").Value = -4067.000100000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 22
$ws.Range("H22").Value = 488.54544
$ws.Range("J22").Value = 433
$ws.Range("L22").Value = 433
$ws.Range("N22").Value = -1133
# row 62
$ws.Range("H62").Value = 2799.6667
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2799.6667
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2799.6667
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4047.6667
# row 65
$ws.Range("H65").Value = 2799.6667
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2799.6667
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 13998.3335
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -20238.3335
# row 86
$ws.Range("H86").Value = 5666.3335
$ws.Range("I86").Value = 5666.3335
$ws.Range("K86").Value = 5666.3335
$ws.Range("M86").Value = -4543.3335
# row 89
$ws.Range("H89").Value = 5666.3335
$ws.Range("I89").Value = 5666.3335
$ws.Range("K89").Value = 28331.6675
$ws.Range("M89").Value = -22715.6675
# row 132
$ws.Range("H132").Value = 2706.5334
$ws.Range("I132").Value = 2471.2856
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 7413.8568
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -4883.8568
$ws.Range("N132").Value = -23060

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 98
$ws.Range("I2").Value = 60.551723
$ws.Range("J2").Value = 149.71428
$ws.Range("K2").Value = 363.310338
$ws.Range("L2").Value = 898.28568
$ws.Range("M2").Value = -250.310338
$ws.Range("N2").Value = -1124.28568
# row 11
$ws.Range("H11").Value = 2369.926
$ws.Range("I11").Value = 2395.6924
$ws.Range("K11").Value = 7187.0772
$ws.Range("M11").Value = -7047.0772
# row 13
$ws.Range("H13").Value = 2045.6
$ws.Range("J13").Value = 3370.3333
$ws.Range("L13").Value = 10110.9999
$ws.Range("N13").Value = -10446.9999
# row 26
$ws.Range("H26").Value = 112.5
$ws.Range("I26").Value = 116.666664
$ws.Range("K26").Value = 349.999992
$ws.Range("M26").Value = -61.99999200000002
# row 47
$ws.Range("H47").Value = 7425.625
$ws.Range("I47").Value = 1343.5714
$ws.Range("K47").Value = 4030.7142
$ws.Range("M47").Value = -3599.7142
# row 51
$ws.Range("H51").Value = 26802.6
$ws.Range("J51").Value = 33005
$ws.Range("L51").Value = 99015
$ws.Range("N51").Value = -99935
# row 112
$ws.Range("H112").Value = 10562.5
$ws.Range("I112").Value = 9928.571
$ws.Range("K112").Value = 29785.713
$ws.Range("M112").Value = -28677.713
# row 114
$ws.Range("H114").Value = 3001
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 3001
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 9003
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -15511
# row 118
$ws.Range("H118").Value = 6588
$ws.Range("I118").Value = 2773.3333
$ws.Range("J118").Value = 8495.333000000001
$ws.Range("K118").Value = 8319.999899999999
$ws.Range("L118").Value = 25485.999
$ws.Range("M118").Value = -7076.999899999999
$ws.Range("N118").Value = -27971.999
# row 122
$ws.Range("H122").Value = 1321.7059
$ws.Range("J122").Value = 1357.4375
$ws.Range("L122").Value = 12216.9375
$ws.Range("N122").Value = -17116.9375
# row 131
$ws.Range("H131").Value = 3035.3845
$ws.Range("J131").Value = 3474.182
$ws.Range("L131").Value = 10422.546
$ws.Range("N131").Value = -20502.546

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 55567344
$ws.Range("I132").Value = 71432020
$ws.Range("K132").Value = 214296060
$ws.Range("M132").Value = -214293530
# row 136
$ws.Range("H136").Value = 40276.5
$ws.Range("J136").Value = 40276.5
$ws.Range("L136").Value = 120829.5
$ws.Range("N136").Value = -125929.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 70237.07000000001
$ws.Range("I7").Value = 2026.8572
$ws.Range("J7").Value = 129921
$ws.Range("K7").Value = 2026.8572
$ws.Range("L7").Value = 129921
$ws.Range("M7").Value = -1914.8572
$ws.Range("N7").Value = -130145
# row 46
$ws.Range("H46").Value = 5201.263
$ws.Range("I46").Value = 1785.4166
$ws.Range("J46").Value = 11057
$ws.Range("K46").Value = 1785.4166
$ws.Range("L46").Value = 11057
$ws.Range("M46").Value = -1597.4166
$ws.Range("N46").Value = -11433
# row 68
$ws.Range("H68").Value = 2050
$ws.Range("J68").Value = 2200
$ws.Range("L68").Value = 2200
$ws.Range("N68").Value = -3698
# row 71
$ws.Range("H71").Value = 2050
$ws.Range("J71").Value = 2200
$ws.Range("L71").Value = 11000
$ws.Range("N71").Value = -18488
# row 126
$ws.Range("H126").Value = 70237.07000000001
$ws.Range("I126").Value = 2026.8572
$ws.Range("J126").Value = 129921
$ws.Range("K126").Value = 6080.571599999999
$ws.Range("L126").Value = 389763
$ws.Range("M126").Value = -3610.571599999999
$ws.Range("N126").Value = -394703
# row 127
$ws.Range("H127").Value = 155950
$ws.Range("J127").Value = 155950
$ws.Range("L127").Value = 155950
$ws.Range("N127").Value = -165870
# row 136
$ws.Range("H136").Value = 245801
$ws.Range("I136").Value = 25000
$ws.Range("J136").Value = 301001.25
$ws.Range("K136").Value = 75000
$ws.Range("L136").Value = 903003.75
$ws.Range("M136").Value = -72450
$ws.Range("N136").Value = -908103.75

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 119
$ws.Range("H119").Value = 90490.5
$ws.Range("J119").Value = 90490.5
$ws.Range("L119").Value = 90490.5
$ws.Range("N119").Value = -100166.5
# row 122
$ws.Range("H122").Value = 18845.607
$ws.Range("I122").Value = 14316.056
$ws.Range("K122").Value = 42948.16800000001
$ws.Range("M122").Value = -40498.16800000001
# row 126
$ws.Range("H126").Value = 9998.714
$ws.Range("J126").Value = 9999.5
$ws.Range("L126").Value = 29998.5
$ws.Range("N126").Value = -34938.5
# row 136
$ws.Range("H136").Value = 1051
$ws.Range("I136").Value = 880.8
$ws.Range("J136").Value = 1334.6666
$ws.Range("K136").Value = 2642.4
$ws.Range("L136").Value = 4003.9998
$ws.Range("M136").Value = -92.39999999999964
$ws.Range("N136").Value = -9103.9998
# row 139
$ws.Range("H139").Value = 76000
$ws.Range("J139").Value = 76000
$ws.Range("L139").Value = 76000
